$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: new columns for curHealth / curDam
$ws.Range("B1").Value = "curHealth"
$ws.Range("C1").Value = "curDam"

# Data rows
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 2

$ws.Range("B3").Value = 5
$ws.Range("C3").Value = 7

$ws.Range("B4").Value = 8
$ws.Range("C4").Value = 8

# Move the active selection to C5, matching the saved selection state
[void]$ws.Range("C5").Select()
